# FedExRateVerification.xlsx - rows 30 & 31 now "pass" verification:
# the ActualRate (col E) is updated to match the ExpectedRate (col D),
# and the Result (col F) flips from FAIL to PASS.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 30: ExpectedRate is $473.23
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '$473.23'
$ws.Range("E30").Style = "Normal"
$ws.Range("F30").Value = "PASS"

# Row 31: ExpectedRate is $252.98
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '$252.98'
$ws.Range("E31").Style = "Normal"
$ws.Range("F31").Value = "PASS"
